# Add the 2022-Q4 quarterly detail sheet and record it on the 总计 (Total)
# summary sheet, pushing the older quarters down one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (Total) summary sheet: insert a new row right under
#    the header for 2022-Q4 and shift the existing quarters down, fixing up
#    the running index in column A as we go.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 15
$total.Cells.Item(2, 4).Value = 4.51
$total.Cells.Item(2, 1).HorizontalAlignment = -4108
$total.Cells.Item(2, 1).VerticalAlignment = -4160
$total.Cells.Item(2, 1).Font.Bold = $true
$total.Cells.Item(2, 1).Borders.LineStyle = 1

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4

# ---------------------------------------------------------------------------
# 2. Create the new "2022-Q4" detail sheet right after "总计" and before
#    the existing "2022-Q1" sheet.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $total)
$ws.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $ws.Cells.Item(1, $c + 2)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @("007128", "天弘增强回报债券A", "44.39", "86.57", "4.45", "1.9754", 6),
    @("007129", "天弘增强回报债券C", "43.01", "86.57", "4.45", "1.9139", 6),
    @("009735", "天弘增强回报债券E", "3.76", "86.57", "4.45", "0.1673", 6),
    @("011899", "长安鑫瑞科技先锋6个月定期开放混合A", "2.72", "81.49", "3.97", "0.1080", 9),
    @("009008", "平安科技创新混合A", "2.49", "92.02", "4.03", "0.1003", 3),
    @("501099", "平安科技创新 3 年封闭混合", "2.73", "91.32", "3.50", "0.0956", 5),
    @("011900", "长安鑫瑞科技先锋6个月定期开放混合C", "1.14", "81.49", "3.97", "0.0453", 9),
    @("009009", "平安科技创新混合C", "0.95", "92.02", "4.03", "0.0383", 3),
    @("014651", "大成专精特新混合A", "0.77", "68.74", "3.87", "0.0298", 10),
    @("740001", "长安宏观策略混合A", "0.26", "90.99", "5.88", "0.0153", 7),
    @("700004", "平安灵活配置混合A", "0.32", "79.07", "2.44", "0.0078", 7),
    @("014652", "大成专精特新混合C", "0.13", "68.74", "3.87", "0.0050", 10),
    @("015078", "平安灵活配置混合C", "0.18", "79.07", "2.44", "0.0044", 7),
    @("016579", "长安宏观策略混合C", "0.03", "90.99", "5.88", "0.0018", 7),
    @("002630", "江信瑞福灵活配置混合A", "0.01", "86.57", "4.45", "0.0004", 6)
)

# Columns B..G carry their numeric-looking values as literal text (to match
# the source data, which keeps things like fund codes / percentages as
# plain strings) -- force the "Text" number format before assignment so the
# engine doesn't silently coerce them into numbers.
$ws.Range("B2:G16").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2

    $idxCell = $ws.Cells.Item($excelRow, 1)
    $idxCell.Value = $r
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1

    $ws.Cells.Item($excelRow, 2).Value = $row[0]
    $ws.Cells.Item($excelRow, 3).Value = $row[1]
    $ws.Cells.Item($excelRow, 4).Value = $row[2]
    $ws.Cells.Item($excelRow, 5).Value = $row[3]
    $ws.Cells.Item($excelRow, 6).Value = $row[4]
    $ws.Cells.Item($excelRow, 7).Value = $row[5]
    $ws.Cells.Item($excelRow, 8).Value = $row[6]
}

$ws.Range("A1").Select()
